$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_3a_Postulate")

# 1. Fix A36: "Z12_B02_P02_Ib01" -> "Z12_B02_P02"
$ws.Range("A36").Value = "Z12_B02_P02"

# 2. Insert a new row before row 48 (shifts old rows 48:50 down to 49:51)
$ws.Rows("48:48").Insert(-4121)  # -4121 = xlShiftDown

# 3. Copy formatting from the row below (now row 49, which holds the old row-48 content)
#    onto the newly inserted blank row 48, so it keeps the same cell style (s="4").
$ws.Range("A49:D49").Copy()
$ws.Range("A48:D48").PasteSpecial(-4122)  # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

# 4. Populate the newly inserted row 48 with the new "Z16_B04" entry
$ws.Range("A48").Value = "Z16_B04_P01"
$ws.Range("B48").Value = "Z16_B04"
$ws.Range("C48").Value = "Negative Auswirkungen auf die politische Ordnung, den Rechtsstaat, die Wirtschaft und die Gesellschaft verhindern"
$ws.Range("D48").Value = "XXXNegative Auswirkungen auf die politische Ordnung, den Rechtsstaat, die Wirtschaft und die Gesellschaft verhindern"
